$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.073.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.859.28'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.19%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.14'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.622'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.95%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.59'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +9.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.329'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.64%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0990'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.127.96'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.42'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.854.38'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.680'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.36%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.096.06'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.30'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.04%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '241.12'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.18'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.70%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.07%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.27'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.94'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +31.78%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.67'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.24%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0557'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.44%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.99'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.90%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.28%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +14.05%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +23.10%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.29'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.27%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.781'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +13.35%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +13.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '91.75'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.01%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.353.68'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.91'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.33'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.20%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +58.67%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.32%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0549'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.34'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.045.37'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0681'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.14%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +16.65%  '

